# Apply updated dSF (column F) values for the richards_garrett workbook.
# This reflects a repull/recalculation of the mean (dSF) data column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -7
    5  = -6
    6  = 6
    7  = -4
    19 = -1
    20 = -3
    23 = -4
    24 = -6
    27 = -2
    33 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
